# Proyecto Ordenes de compra
# Adds a new project entry ("Seguimiento de Ocs" / Proy_TCon_005) as a new
# row in the "Tabla1" table on Hoja1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$lo = $ws.ListObjects.Item("Tabla1")

# Add a new row to the table - this grows the table range (and the
# autoFilter range) from A1:E10 to A1:E11 automatically.
$newRow = $lo.ListRows.Add()

# The newly inserted row does not automatically inherit the banded-row /
# bordered formatting used by the rest of the table, so copy the format
# from the previous last data row (row 10) onto the new row (row 11).
$ws.Range("A10:E10").Copy()
$newRow.Range.PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Fill in the new record's values.
$newRow.Range.Cells.Item(1, 1).Value = "Seguimiento de Ocs"
$newRow.Range.Cells.Item(1, 2).Value = "Proy_TCon_005"
$newRow.Range.Cells.Item(1, 3).Value = "portada_proyecto_010.png"
$newRow.Range.Cells.Item(1, 4).Value = "https://app.powerbi.com/view?r=eyJrIjoiYzUyMDdmNmUtY2Y1Yy00YzY2LWJmZDQtZGIyMTQyYjQ1OTJmIiwidCI6IjkxZjVjYjg5LTUyZmUtNDdhYi05MDVmLTRlMzU4ODZmNWE1NyIsImMiOjR9"

# Match the author's final selection / active cell on the sheet.
$ws.Range("A11").Select()
